$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.019.76"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "2.921.90"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.44"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.04"
$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.506"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "2.921.05"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.85"
$ws.Range("E10").Value = "  +3.10%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("E13").Value = "  +1.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.64"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "3.406.29"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "61.019.06"
$ws.Range("E17").Value = "  +0.94%  "

$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "2.918.60"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.00"
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.09"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.03"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  +3.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  +2.29%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +5.98%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.12"
$ws.Range("E32").Value = "  -1.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.53"
$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("D35").Value = "0.0₃0861"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +4.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.97"
$ws.Range("E39").Value = "  +0.66%  "

$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.41"
$ws.Range("E44").Value = "  -4.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "376.15"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0348"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").Value = "2.710.79"
$ws.Range("E47").Value = "  +2.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.66"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.19"
$ws.Range("E50").Value = "  -3.80%  "

$ws.Range("E51").Value = "  +0.50%  "
